$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry data which differs from row to row.
$cols = @("A","B","D","E","F","G","H","I","J","Q","R","AC")

# Mapping: for target row (key), the row whose CURRENT (pre-edit) data
# must be copied into it (value). This reproduces the row permutation
# described by the diff.
$srcOf = @{
    2  = 12
    3  = 13
    4  = 2
    5  = 3
    6  = 4
    7  = 5
    8  = 14
    9  = 15
    10 = 16
    11 = 6
    12 = 7
    13 = 8
    14 = 9
    15 = 10
    16 = 11
}

# 1) Snapshot all current values for the relevant columns/rows before
#    writing anything (rows are being rearranged amongst each other).
$snapshot = @{}
foreach ($r in 2..16) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range($c + $r).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write back the permuted values.
foreach ($r in 2..16) {
    $src = $srcOf[$r]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $addr = $c + $r
        $val = $rowVals[$c]
        if ($null -eq $val -or $val -eq "") {
            $ws.Range($addr).ClearContents()
        }
        elseif ($c -eq "I" -or $c -eq "J") {
            # These columns hold numeric-looking text (e.g. "2", "20").
            # Prefix with an apostrophe so Excel keeps them as text
            # instead of silently converting them to numbers.
            $ws.Range($addr).Value = "'" + $val
        }
        else {
            $ws.Range($addr).Value = $val
        }
    }
}
